$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-01-04 Sunday" "2026-01-05 Monday"

Replace-Text "823×4=" "746×4="
Replace-Text "627×8=" "653×2="
Replace-Text "796×5=" "177×7="
Replace-Text "285×3=" "624×9="
Replace-Text "998×7=" "806×3="

Replace-Text "212×3=" "200×7="
Replace-Text "972×5=" "609×3="
Replace-Text "855×8=" "885×5="
Replace-Text "911×5=" "927×2="
Replace-Text "271×8=" "355×6="

Replace-Text "733×7=" "559×8="
Replace-Text "339×8=" "937×8="
Replace-Text "560×8=" "209×9="
Replace-Text "179×3=" "618×6="
Replace-Text "725×6=" "795×2="

Replace-Text "462×4=" "578×2="
Replace-Text "323×4=" "924×6="
Replace-Text "361×9=" "180×3="
Replace-Text "417×5=" "878×7="
Replace-Text "363×5=" "262×4="

Replace-Text "852×8=" "988×2="
Replace-Text "144×7=" "483×4="
Replace-Text "819×6=" "224×3="
Replace-Text "912×8=" "365×3="
Replace-Text "243×2=" "331×8="
